# This script re-assigns the "unit" label (column C) of each of the five
# sections on the sheet so that every unit gets its own distinct label
# instead of the previous sequential unit1..unit5 pattern that repeated
# across sections (see commit message: "refactor section unit day logic
# to assign individual days to each unit rather than a sequential pattern
# across sections").
#
# Mapping discovered from the target diff (old label text -> new label text),
# applied per physical row-range/section rather than simply by old text,
# because "unit5" appears in two different sections that must resolve to
# two different new labels:
#   rows   4-87   : unit1        -> unit1   (unchanged)
#   rows 102-174  : unit2        -> unit1
#   rows 213-299  : unit3        -> unit2
#   rows 332-403  : unit4        -> unit2
#   rows 441-508  : unit5        -> unit3
#   row  442 only : " unit5"     -> " unit3"   (stray leading-space variant)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sections = @(
    @{ Start = 4;   End = 87;  New = "unit1" },
    @{ Start = 102; End = 174; New = "unit1" },
    @{ Start = 213; End = 299; New = "unit2" },
    @{ Start = 332; End = 403; New = "unit2" },
    @{ Start = 441; End = 508; New = "unit3" }
)

foreach ($section in $sections) {
    for ($r = $section.Start; $r -le $section.End; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        $cur = $cell.Value2
        if ($cur -eq $null -or $cur -eq "") {
            continue
        }
        if ($cur -eq " unit5") {
            $cell.Value = " unit3"
        } else {
            $cell.Value = $section.New
        }
    }
}

# Update the saved view state to match the new scroll/selection position.
$ws.Range("H51").Select()
